# Update Name of Algo
# Apply updated imputed values to result_data_KNN sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.336
$ws.Range("A3").Value = -21.557
$ws.Range("B5").Value = 6.712999999999999
$ws.Range("D5").Value = -8.440999999999999
$ws.Range("D9").Value = -7.904999999999998
$ws.Range("D11").Value = -8.359
$ws.Range("A14").Value = -20.891
$ws.Range("A16").Value = -20.947
$ws.Range("B16").Value = 6.425
$ws.Range("D17").Value = -7.920999999999999
$ws.Range("A21").Value = -21.284
$ws.Range("D21").Value = -7.825999999999999
$ws.Range("A23").Value = -21.666
$ws.Range("A25").Value = -22.078
